$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "hasTimeStamp" property row (row 6): the file-specific
# timestamp labels (B:E) gain a "(file)" / "(Datei)" / "(fichier)" suffix
# to disambiguate them from the generic "Time stamp" labels in G:J.
$ws.Range("B6").Value = "Time Stamp (file)"
$ws.Range("C6").Value = "Zeitstempel (Datei)"
$ws.Range("D6").Value = "Horodatage (fichier)"
$ws.Range("E6").Value = "Timestamp (file)"
# Touch the font so Excel records an explicit (if visually identical)
# font/style entry for these re-typed cells, matching a manual edit.
$ws.Range("B6:E6").Font.Name = "Calibri"

# --- Add a new ontology property row: "hasSeqnum" (row 10), describing a
# sequence number used for compound objects, in all supported languages.
$ws.Range("A10").Value = "hasSeqnum"
$ws.Range("B10").Value = "Seqnum"
$ws.Range("C10").Value = "Seqnum"
$ws.Range("D10").Value = "Seqnum"
$ws.Range("E10").Value = "Seqnum"
$ws.Range("G10").Value = "Sequence number used for compound object"
$ws.Range("H10").Value = "Sequenznummer für zusammengesetzte Objekte"
$ws.Range("I10").Value = "Numéro de séquence utilisé pour l'objet composé"
$ws.Range("J10").Value = "Numero di sequenza utilizzato per oggetti composti"
$ws.Range("L10").Value = "seqnum, schema:position"
$ws.Range("M10").Value = "IntValue"
$ws.Range("N10").Value = "SimpleText"

# Give the new row its own (visually default) explicit font, as Excel does
# for freshly entered data, and mark the "super" cell slightly differently
# as happened in the source edit.
$ws.Range("A10:J10").Font.Name = "Calibri"
$ws.Range("N10").Font.Name = "Calibri"
$ws.Range("L10").Font.Name = "Calibri"

# --- Widen column K (super/comment column) to fit the longer new content.
$ws.Range("K:K").ColumnWidth = 39.67

# --- Move the active selection, as left by the editor after the change.
$ws.Range("M14").Select()
